# Insert a new weekly record as row 50, pushing existing rows 50-70 down to 51-71.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(50).Insert()

$ws.Cells.Item(50, 1).Value  = 1
$ws.Cells.Item(50, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(50, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(50, 4).Value  = 44489
$ws.Cells.Item(50, 5).Value  = 15
$ws.Cells.Item(50, 6).Value  = "Fruta"
$ws.Cells.Item(50, 7).Value  = 100102
$ws.Cells.Item(50, 8).Value  = "Cítricos"
$ws.Cells.Item(50, 9).Value  = 100102005
$ws.Cells.Item(50, 10).Value = "Naranja"
$ws.Cells.Item(50, 11).Value = "Lane Late"
$ws.Cells.Item(50, 12).Value = "Segunda"
$ws.Cells.Item(50, 13).Value = 250
$ws.Cells.Item(50, 14).Value = 650
$ws.Cells.Item(50, 15).Value = 700
$ws.Cells.Item(50, 16).Value = 675
$ws.Cells.Item(50, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(50, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(50, 19).Value = 675
$ws.Cells.Item(50, 20).Value = 1
